# Applies the re-ordering + appended-matches update described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rA = $ws.Range("F$($rowA):V$($rowA)")
    $rB = $ws.Range("F$($rowB):V$($rowB)")
    $vA = $rA.Value2
    $vB = $rB.Value2
    $rA.Value = $vB
    $rB.Value = $vA
}

# Rows 89 <-> 90 swap (same kickoff date, re-ordered)
Swap-Rows 89 90

# Rows 94 <-> 95 swap
Swap-Rows 94 95

# Rows 100 <-> 101 swap
Swap-Rows 100 101

# Rows 102 <-> 103 swap
Swap-Rows 102 103

# Rows 125-128 rotate: 125<-126, 126<-127, 127<-128, 128<-125
$r125 = $ws.Range("F125:V125")
$r126 = $ws.Range("F126:V126")
$r127 = $ws.Range("F127:V127")
$r128 = $ws.Range("F128:V128")

$v125 = $r125.Value2
$v126 = $r126.Value2
$v127 = $r127.Value2
$v128 = $r128.Value2

$r125.Value = $v126
$r126.Value = $v127
$r127.Value = $v128
$r128.Value = $v125

# Append two new rows (129, 130) with the same formatting as the row above.
$ws.Range("A128:V128").Copy($ws.Range("A129:V129"))
$ws.Range("A128:V128").Copy($ws.Range("A130:V130"))

# Row 129 data
$ws.Range("A129").Value = 128
$ws.Range("B129").Value = "poland"
$ws.Range("C129").Value = "iii-liga-group-ii"
$ws.Range("D129").Value = "2023-2024"
$ws.Range("E129").Value = 45235.54166666666
$ws.Range("F129").Value = "Cartusia Kartuzy"
$ws.Range("G129").Value = 1
$ws.Range("H129").Value = "Swinoujscie"
$ws.Range("I129").Value = 2
$ws.Range("J129").Value = 1.49
$ws.Range("K129").Value = "04/11/2023 01:13"
$ws.Range("L129").Value = 1.55
$ws.Range("M129").Value = "05/11/2023 12:51"
$ws.Range("N129").Value = 4.04
$ws.Range("O129").Value = "04/11/2023 01:13"
$ws.Range("P129").Value = 4.06
$ws.Range("Q129").Value = "05/11/2023 12:51"
$ws.Range("R129").Value = 4.33
$ws.Range("S129").Value = "04/11/2023 01:13"
$ws.Range("T129").Value = 4.62
$ws.Range("U129").Value = "05/11/2023 12:51"
$ws.Range("V129").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-ii/cartusia-kartuzy-swinoujscie/KCgjQZJe/"

# Row 130 data
$ws.Range("A130").Value = 129
$ws.Range("B130").Value = "poland"
$ws.Range("C130").Value = "iii-liga-group-ii"
$ws.Range("D130").Value = "2023-2024"
$ws.Range("E130").Value = 45235.54166666666
$ws.Range("F130").Value = "Zawisza"
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = "Swit Skolwin"
$ws.Range("I130").Value = 1
$ws.Range("J130").Value = 2.55
$ws.Range("K130").Value = "04/11/2023 01:13"
$ws.Range("L130").Value = 2.69
$ws.Range("M130").Value = "05/11/2023 12:05"
$ws.Range("N130").Value = 3.19
$ws.Range("O130").Value = "04/11/2023 01:13"
$ws.Range("P130").Value = 3.16
$ws.Range("Q130").Value = "05/11/2023 12:05"
$ws.Range("R130").Value = 2.27
$ws.Range("S130").Value = "04/11/2023 01:13"
$ws.Range("T130").Value = 2.4
$ws.Range("U130").Value = "05/11/2023 12:05"
$ws.Range("V130").Value = "https://www.betexplorer.com/football/poland/iii-liga-group-ii/zawisza-swit-skolwin/GbwgpeBR/"

